$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.140.43'
$ws.Range('D3').Value = '2.307.13'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'301.78"
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').Value = "'99.50"
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('D7').Value = "'0.505"
$ws.Range('E7').Value = '  +2.37%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.506"
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').Value = "'34.31"
$ws.Range('E10').Value = '  +4.05%  '
$ws.Range('D11').Value = "'0.0798"
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('E12').Value = '  +4.12%  '
$ws.Range('D13').Value = "'17.93"
$ws.Range('E13').Value = '  +14.97%  '
$ws.Range('D14').Value = "'6.80"
$ws.Range('E14').Value = '  +2.70%  '
$ws.Range('D15').Value = '2.665.52'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '2.319.84'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = "'0.814"
$ws.Range('E17').Value = '  +5.57%  '
$ws.Range('D18').Value = '43.048.36'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = "'12.59"
$ws.Range('E19').Value = '  +11.65%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').Value = "'6.11"
$ws.Range('E21').Value = '  +2.67%  '
$ws.Range('D22').Value = "'67.73"
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').Value = "'237.07"
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').Value = "'2.22"
$ws.Range('E24').Value = '  +15.08%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').Value = "'24.76"
$ws.Range('E27').Value = '  +4.24%  '
$ws.Range('D28').Value = "'168.65"
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  -8.59%  '
$ws.Range('D30').Value = "'34.11"
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').Value = "'9.15"
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').Value = "'1.00"
$ws.Range('D33').Value = "'5.05"
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').Value = "'2.43"
$ws.Range('E34').Value = '  +4.39%  '
$ws.Range('D35').Value = "'4.54"
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').Value = "'17.04"
$ws.Range('E36').Value = '  +6.39%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  +4.05%  '
$ws.Range('E39').Value = '  +5.59%  '
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('D42').Value = "'2.36"
$ws.Range('E42').Value = '  -4.76%  '
$ws.Range('D43').Value = '1.991.93'
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('E45').Value = '  +5.60%  '
$ws.Range('D46').Value = "'17.73"
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('E47').Value = '  +2.90%  '
$ws.Range('D48').Value = "'56.44"
$ws.Range('E48').Value = '  +9.02%  '
$ws.Range('D49').Value = '2.534.24'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('E50').Value = '  +5.09%  '
$ws.Range('E51').Value = '  +1.25%  '
